# "added photo to excel option" -- populate the new "Pax" (column C) values
# for rows 2-31 on Sheet1, and leave the sheet scrolled to the top with
# C2:C31 selected (matching the author's final view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$paxValues = @(10, 10, 10, 10, 10, 10, 10, 10, 10, 10, 10, 10, 12, 10, 10, 10, 10, 10, 10, 10, 10, 10, 8, 10, 10, 8, 8, 8, 8, 8)

for ($i = 0; $i -lt $paxValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $paxValues[$i]
}

# Scroll the window back to the top of the sheet and select C2:C31, like the
# saved workbook's final view.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C2:C31").Select()
